$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test method names to reflect refactor (Purchase/Prorated/Upgrade -> Downgrade wording)
$ws.Range("A4").Value = "verifyCustomerAdditionalLicenseDowngrade"
$ws.Range("A6").Value = "verifyCustomerReceiptPageWithRecurringOrderDetails"
$ws.Range("A7").Value = "verifyCustomerReceivedSubscriptionDowngradeReceipt"

# Adjust column A width to fit the now-shorter text
# (ColumnWidth is in "characters"; OOXML stored width = characters + 0.8333333333333334,
#  so back-solve for the characters value that yields the target stored width of 50.9921875)
$ws.Columns.Item(1).ColumnWidth = 50.158854166666664
